$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5132.7827
$ws.Range("I40").Value = 4610.375
$ws.Range("K40").Value = 4610.375
$ws.Range("M40").Value = -4435.375
$ws.Range("H74").Value = 18934.05
$ws.Range("I74").Value = 18763.883
$ws.Range("J74").Value = 19898.334
$ws.Range("K74").Value = 18763.883
$ws.Range("L74").Value = 19898.334
$ws.Range("M74").Value = -17827.883
$ws.Range("N74").Value = -21770.334
$ws.Range("H77").Value = 18934.05
$ws.Range("I77").Value = 18763.883
$ws.Range("J77").Value = 19898.334
$ws.Range("K77").Value = 93819.41500000001
$ws.Range("L77").Value = 99491.67
$ws.Range("M77").Value = -89139.41500000001
$ws.Range("N77").Value = -108851.67
$ws.Range("H100").Value = 1986.1666
$ws.Range("I100").Value = 2103.4
$ws.Range("J100").Value = 1400
$ws.Range("K100").Value = 2103.4
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -1562.4
$ws.Range("N100").Value = -2482
$ws.Range("H112").Value = 2840.1667
$ws.Range("J112").Value = 2840.1667
$ws.Range("L112").Value = 8520.500100000001
$ws.Range("N112").Value = -10736.5001
$ws.Range("H118").Value = 1133
$ws.Range("I118").Value = 1216
$ws.Range("J118").Value = 1050
$ws.Range("K118").Value = 3648
$ws.Range("L118").Value = 3150
$ws.Range("M118").Value = -1991
$ws.Range("N118").Value = -6464
$ws.Range("H132").Value = 9872.0625
$ws.Range("I132").Value = 9872.0625
$ws.Range("K132").Value = 29616.1875
$ws.Range("M132").Value = -27086.1875
$ws.Range("H137").Value = 3324.6155
$ws.Range("I137").Value = 1371.6
$ws.Range("J137").Value = 5987.8184
$ws.Range("K137").Value = 4114.799999999999
$ws.Range("L137").Value = 17963.4552
$ws.Range("M137").Value = -1564.799999999999
$ws.Range("N137").Value = -23063.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 609.46875
$ws.Range("I2").Value = 604.53845
$ws.Range("J2").Value = 630.8333
$ws.Range("K2").Value = 604.53845
$ws.Range("L2").Value = 630.8333
$ws.Range("M2").Value = -491.53845
$ws.Range("N2").Value = -856.8333
$ws.Range("H116").Value = 609.46875
$ws.Range("I116").Value = 604.53845
$ws.Range("J116").Value = 630.8333
$ws.Range("K116").Value = 604.53845
$ws.Range("L116").Value = 630.8333
$ws.Range("M116").Value = 1689.46155
$ws.Range("N116").Value = -5218.8333
$ws.Range("H122").Value = 5704.25
$ws.Range("I122").Value = 5328.5625
$ws.Range("K122").Value = 15985.6875
$ws.Range("M122").Value = -13535.6875
$ws.Range("H139").Value = 70962.5
$ws.Range("J139").Value = 70962.5
$ws.Range("L139").Value = 70962.5
$ws.Range("N139").Value = -81242.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 609.46875
$ws.Range("I3").Value = 604.53845
$ws.Range("J3").Value = 630.8333
$ws.Range("K3").Value = 604.53845
$ws.Range("L3").Value = 630.8333
$ws.Range("M3").Value = -490.53845
$ws.Range("N3").Value = -858.8333
$ws.Range("H20").Value = 25257038
$ws.Range("I20").Value = 30868886
$ws.Range("K20").Value = 30868886
$ws.Range("M20").Value = -30868639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1728.8
$ws.Range("J16").Value = 575
$ws.Range("L16").Value = 575
$ws.Range("N16").Value = -1149
$ws.Range("H31").Value = 3941.5881
$ws.Range("I31").Value = 3563.3044
$ws.Range("J31").Value = 4252.3213
$ws.Range("K31").Value = 3563.3044
$ws.Range("L31").Value = 4252.3213
$ws.Range("M31").Value = -3268.3044
$ws.Range("N31").Value = -4842.3213
$ws.Range("H34").Value = 3941.5881
$ws.Range("I34").Value = 3563.3044
$ws.Range("J34").Value = 4252.3213
$ws.Range("K34").Value = 3563.3044
$ws.Range("L34").Value = 4252.3213
$ws.Range("M34").Value = -3361.3044
$ws.Range("N34").Value = -4656.3213
$ws.Range("H58").Value = 4452.143
$ws.Range("J58").Value = 4529.7144
$ws.Range("L58").Value = 4529.7144
$ws.Range("N58").Value = -4935.7144
$ws.Range("H113").Value = 1728.8
$ws.Range("J113").Value = 575
$ws.Range("L113").Value = 575
$ws.Range("N113").Value = -4915
$ws.Range("H134").Value = 3205.158
$ws.Range("I134").Value = 2935.0715
$ws.Range("J134").Value = 3961.4
$ws.Range("K134").Value = 8805.2145
$ws.Range("L134").Value = 11884.2
$ws.Range("M134").Value = -6270.2145
$ws.Range("N134").Value = -16954.2
$ws.Range("H136").Value = 4452.143
$ws.Range("J136").Value = 4529.7144
$ws.Range("L136").Value = 13589.1432
$ws.Range("N136").Value = -18689.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1000428
$ws.Range("I97").Value = 2500179
$ws.Range("K97").Value = 7500537
$ws.Range("M97").Value = -7500041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 997
$ws.Range("I107").Value = 997
$ws.Range("K107").Value = 997
$ws.Range("M107").Value = 923
$ws.Range("H113").Value = 5458.591
$ws.Range("I113").Value = 5356.6313
$ws.Range("K113").Value = 5356.6313
$ws.Range("M113").Value = -3186.6313
$ws.Range("H122").Value = 2468.5557
$ws.Range("I122").Value = 2888.2856
$ws.Range("K122").Value = 8664.856800000001
$ws.Range("M122").Value = -6214.856800000001
$ws.Range("H132").Value = 2430.074
$ws.Range("I132").Value = 2064.48
$ws.Range("K132").Value = 6193.440000000001
$ws.Range("M132").Value = -3663.440000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 69998.5
$ws.Range("J87").Value = 69998.5
$ws.Range("L87").Value = 69998.5
$ws.Range("N87").Value = -72244.5
$ws.Range("H90").Value = 69998.5
$ws.Range("J90").Value = 69998.5
$ws.Range("L90").Value = 209995.5
$ws.Range("N90").Value = -221227.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22731038
$ws.Range("I122").Value = 4181.125
$ws.Range("K122").Value = 12543.375
$ws.Range("M122").Value = -10093.375
$ws.Range("H126").Value = 2886
$ws.Range("I126").Value = 1477.1666
$ws.Range("J126").Value = 4999.25
$ws.Range("K126").Value = 4431.4998
$ws.Range("L126").Value = 14997.75
$ws.Range("M126").Value = -1961.4998
$ws.Range("N126").Value = -19937.75
$ws.Range("H132").Value = 7754353
$ws.Range("I132").Value = 8774284
$ws.Range("J132").Value = 2880
$ws.Range("K132").Value = 26322852
$ws.Range("L132").Value = 8640
$ws.Range("M132").Value = -26320322
$ws.Range("N132").Value = -13700

Write-Output "Applied all updates"